# Updates cryptos list values (price & 1h volume change) per the
# "Updated cryptos list on Sun Feb 18 11:12:12 UTC 2024 with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.797.55"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.807.96"
$ws.Range("E3").Value = "  +1.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.10"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.49"
$ws.Range("E6").Value = "  +1.99%  "

# Row 7
$ws.Range("E7").Value = "  +0.72%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +8.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.47"
$ws.Range("E10").Value = "  +2.27%  "

# Row 11
$ws.Range("E11").Value = "  -0.89%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("E12").Value = "  -0.47%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.06"
$ws.Range("E13").Value = "  +3.02%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.79"
$ws.Range("E14").Value = "  +2.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.251.76"
$ws.Range("E15").Value = "  +1.12%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.824.03"
$ws.Range("E16").Value = "  +1.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.946"
$ws.Range("E17").Value = "  +1.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.784.30"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("E19").Value = "  +2.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.20"
$ws.Range("E20").Value = "  +3.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.66"
$ws.Range("E21").Value = "  +4.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0978"
$ws.Range("E22").Value = "  +1.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.58"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.19"
$ws.Range("E24").Value = "  +0.53%  "

# Row 25
$ws.Range("E25").Value = "  +1.56%  "

# Row 26
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("E27").Value = "  -0.55%  "

# Row 28
$ws.Range("E28").Value = "  -1.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.60"
$ws.Range("E29").Value = "  +11.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.41"
$ws.Range("E30").Value = "  +2.46%  "

# Row 31
$ws.Range("E31").Value = "  +3.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.44"
$ws.Range("E32").Value = "  +1.95%  "

# Row 33
$ws.Range("E33").Value = "  +0.20%  "

# Row 34
$ws.Range("E34").Value = "  +9.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0889"
$ws.Range("E35").Value = "  +6.44%  "

# Row 36
$ws.Range("E36").Value = "  -0.35%  "

# Row 37
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("E38").Value = "  -0.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("E39").Value = "  +2.65%  "

# Row 40
$ws.Range("E40").Value = "  +0.97%  "

# Row 41
$ws.Range("E41").Value = "  +1.40%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  -1.79%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.41"
$ws.Range("E43").Value = "  +0.66%  "

# Row 44
$ws.Range("E44").Value = "  -0.19%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.08"
$ws.Range("E45").Value = "  +1.92%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.43"
$ws.Range("E46").Value = "  +5.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.111.90"
$ws.Range("E47").Value = "  +0.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +6.77%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.952"
$ws.Range("E49").Value = "  +0.80%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.48"
$ws.Range("E50").Value = "  -1.02%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.37"
$ws.Range("E51").Value = "  +6.81%  "
